$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet (part 1): insert four new explanatory note rows (18-21)
# plus a blank spacer row (22) before the existing LCFS note, which shifts
# from rows 18-20 down to rows 23-25.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Rows.Item(18).Insert()
$about.Rows.Item(18).Insert()
$about.Rows.Item(18).Insert()
$about.Rows.Item(18).Insert()
$about.Rows.Item(18).Insert()
# Inserting copies the (bold) format from the row above (row 17, "Notes:");
# put the new block back to normal (unbolded) formatting.
$about.Range("A18:A22").Font.Bold = $false

$about.Range("A18").Value = "For vehicle types that can use electricity, this variable specifies the percentage"
$about.Range("A19").Value = "reduction in fuel use (on a BTU basis) relative to the typical fuel type for that vehicle"
$about.Range("A20").Value = "type (e.g. gasoline for LDVs, diesel for HDVs, etc.) due to the fact that electricity"
$about.Range("A21").Value = "can be converted into work more efficiently than other fuel types."
# row 22 intentionally stays blank (spacer before the pre-existing LCFS note)

Write-Output "about sheet notes inserted"

# ---------------------------------------------------------------------------
# "PTFURfE" sheet: relabel the header, widen column A, taller header row,
# and flow the "aircraft"/"ships" rows (4 & 6) from hard-coded zeros to
# formulas mirroring the "rail" row (5).
# ---------------------------------------------------------------------------
$ptf = $wb.Worksheets.Item("PTFURfE")

$ptf.Range("A1").Value = "Percentage Reduction (dimensionless)"
$ptf.Range("A1").Font.Bold = $true
$ptf.Range("A1").WrapText = $true
$ptf.Rows.Item(1).RowHeight = 45

$ptf.Columns.Item(1).ColumnWidth = 16.85546875

$ptf.Range("B4").Formula = "=B5"
$ptf.Range("C4").Formula = "=C5"
$ptf.Range("B4:C4").NumberFormat = $ptf.Range("B5").NumberFormat

$ptf.Range("B6").Formula = "=B5"
$ptf.Range("C6").Formula = "=C5"
$ptf.Range("B6:C6").NumberFormat = $ptf.Range("B5").NumberFormat

Write-Output "PTFURfE sheet updated"

# ---------------------------------------------------------------------------
# "About" sheet (part 2): append a new two-line note at the end (rows
# 27-28, leaving row 26 blank).
# ---------------------------------------------------------------------------
$about.Range("A27").Value = "Aircraft and ships are assumed to be the same as rail, since they all use large engines"
$about.Range("A28").Value = "intended to move heavy craft."

Write-Output "about sheet closing note appended"
